$p = $ppt.ActivePresentation
$s = $p.Slides.Item(29)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Original text in this shape begins: "it('should have one child', () => {    expect(...
# Change "one child'" -> "two children'" inside the bold/green quoted string; this
# splits the existing bold/green run into two runs that both keep its formatting
# (same rPr/fill), matching 'should have ' + 'two children'.
$quoted = $tr.Characters(17, 10)
$quoted.Text = "two children'"

# After the above edit, the following ", () => {" run has shifted right by 3
# characters (because "two children'" is 3 characters longer than "one child'").
# Re-split it into its own ", " run (formatting unchanged) followed by the
# untouched "() => {" run, matching the target run layout.
$tail = $tr.Characters(30, 2)
$tail.Text = $tail.Text
